$wb = $excel.ActiveWorkbook

# --- Workbook-level change ---
# The defined name "cgam_processes" now covers columns A:D instead of A:E
# (Processes sheet lost a column from the named range).
$wb.Names.Item("cgam_processes").RefersTo = "=Processes!`$A`$1:`$D`$7"

# --- Processes sheet changes ---
$ws = $wb.Worksheets.Item("Processes")

# The "description" and "type" columns were swapped: what used to be column B
# (description) is now column E, and what used to be column E (type) is now
# column B. Swap the cell contents for the header row and the 11 data rows.
for ($r = 1; $r -le 12; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $eCell = $ws.Cells.Item($r, 5)
    $bVal = $bCell.Value2
    $eVal = $eCell.Value2
    $bCell.Value = $eVal
    $eCell.Value = $bVal
}

# Column widths follow the new content: B (now "type") is narrower,
# E (now "description") is wider.
$ws.Columns.Item(2).ColumnWidth = 13.29
$ws.Columns.Item(5).ColumnWidth = 33.43

# The process-type list validation moves from the old column E to the new
# column B (it stays attached to the "type" data, which now lives in B).
$rng = $ws.Range("B2:B12")
$rng.Validation.Add(3, 2, 1, "=Validate!`$B`$2:`$B`$3")
$rng.Validation.ErrorTitle = "Process Type"
$rng.Validation.ErrorMessage = "Invalid Process Type"
$rng.Validation.IgnoreBlank = $true
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $true
$rng.Validation.ShowError = $true

# Reflect the new selection/active cell on the sheet (now column E, the
# "description" column, fully selected).
[void]$ws.Activate()
$ws.Range("E1:E12").Select() | Out-Null
